$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-Alert Template Import")

# Constant for xlPasteFormats
$xlPasteFormats = -4122

# Scratch cell far outside the used range (A1:AP3), used to stash/restore
# a cell's direct formatting (including attributes like quotePrefix that
# are not reachable through simple COM properties) across a value change.
$scratch = $ws.Range("ZZ100")

function Set-TextPreserveFormat {
    param(
        [string]$addr,
        [string]$newValue
    )
    $cell = $ws.Range($addr)

    # Stash current formatting onto the scratch cell.
    $cell.Copy() | Out-Null
    $scratch.PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = 0

    # Update the value (this may reset the cell's direct formatting).
    $cell.Value = $newValue

    # Restore the original formatting from the scratch cell.
    $scratch.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = 0
}

Set-TextPreserveFormat "A3" "JSSO1000246"
Set-TextPreserveFormat "B3" "JSSO1000246"
Set-TextPreserveFormat "C3" "JSSO1000246"
Set-TextPreserveFormat "AJ3" "JSCN1000246"
Set-TextPreserveFormat "AN3" "MBLJSSO1000246"
Set-TextPreserveFormat "AO3" "HBLJSSO1000246"

# Clean up the scratch cell so no stray formatting/content is left behind.
$scratch.Clear() | Out-Null
